$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data: Date (46006) and Error Count (93)
$ws.Range("A34").Value = 46006
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B34").Value = 93

# Update the selection to mirror the recorded view state
$ws.Range("A34:B34").Select()
